# Backup QR Scanner data - 2025-12-29T11:57:03.079Z - Cache Bust: 1767009423079
#
# 1) Rename the worksheet tab from "Session" to "Neurology"
# 2) Append a new scan-log row (row 86) with the latest "Manual" entry,
#    keeping every column stored as text (matching the existing rows),
#    which also pushes the sheet's dimension out to A1:F86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Neurology"

# New row 86: a manually-logged scan event. Every value on this sheet is
# kept as literal text (student IDs, dates and times included), so each
# cell is written as a formula whose result is the quoted string literal
# — this is what keeps the stored type as text (not a number/date).
$ws.Range("A86").Formula = "=""201560"""
$ws.Range("B86").Formula = "=""Neurology"""
$ws.Range("C86").Formula = "=""29/12/2025"""
$ws.Range("D86").Formula = "=""13:57:00"""
$ws.Range("E86").Formula = "=""Manual"""
$ws.Range("F86").Formula = "=""emp17.farah.a.youssef@gmail.com"""
